# Stock update / stock create flows against the ProductData sheet.
#
# Row 7 (existing item "1P") gets its quantity/purchase_price/sale_price
# normalised from text to real numbers as part of StockManageView's update
# (the original row had those three columns stored as text).
#
# Two new rows are appended:
#   - Row 8: a brand-new stock item created through the "update" path,
#     which stores quantity/prices as real numbers.
#   - Row 9: a brand-new stock item created through the plain "create" path,
#     which (like the legacy row 7 before it) stores quantity/prices as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up row 7: cast quantity / purchase_price / sale_price to numbers ---
$ws.Range("D7").Value = 45
$ws.Range("E7").Value = 1200
$ws.Range("F7").Value = 2500

# Reuse the existing creation_date date/time format for any new timestamp cell.
$dateFormat = $ws.Range("G7").NumberFormat()

# --- Row 8: new stock item, numeric quantity/prices ---
$ws.Cells.Item(8, 1).Value = "1L"
$ws.Cells.Item(8, 2).Value = "Limpieza"
$ws.Cells.Item(8, 3).Value = "Limpido Clorox"
$ws.Cells.Item(8, 4).Value = 24
$ws.Cells.Item(8, 5).Value = 12000
$ws.Cells.Item(8, 6).Value = 16000

# Stamp the creation date with "now" (evaluated once, then frozen to a plain
# value so the saved cell holds a literal timestamp, not a live formula).
$ws.Range("G8").Formula = "=NOW()"
$g8Now = $ws.Range("G8").Value()
$ws.Range("G8").Value = $g8Now
$ws.Range("G8").NumberFormat = $dateFormat

# --- Row 9: new stock item, quantity/prices kept as plain text ---
$ws.Cells.Item(9, 1).Value = "1AP"
$ws.Cells.Item(9, 2).Value = "Aseo personal"
$ws.Cells.Item(9, 3).Value = "Jabon piel dove en barra x 4"

foreach ($addr in @("D9", "E9", "F9")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("D9").Value = "30"
$ws.Range("E9").Value = "5600"
$ws.Range("F9").Value = "8500"
foreach ($addr in @("D9", "E9", "F9")) {
    $ws.Range($addr).Style = "Normal"
}

$ws.Range("G9").Formula = "=NOW()"
$g9Now = $ws.Range("G9").Value()
$ws.Range("G9").Value = $g9Now
$ws.Range("G9").NumberFormat = $dateFormat
